# Updates the crypto price/volume table (cols D & E) for each coin row,
# and swaps the XRP / USDC rows (7 and 8) to their new order, matching the
# "Updated cryptos list" GitHub Actions commit.
#
# Cell values like "1.00", "191.96", "0.0690" etc. look numeric, so a plain
# Range.Value assignment would make Excel auto-convert them to a number
# (dropping the trailing zero / decimal formatting the source data relies
# on). Set-TextValue forces those through as literal text while leaving the
# cell's style untouched (format as text only for the assignment, then
# ClearFormats so no stray number-format / quote-prefix style sticks to the
# cell afterwards).
function Set-TextValue {
    param($Worksheet, [string]$CellRef, [string]$Text)
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.531.18'
$ws.Range("E2").Value = '  +3.13%  '
$ws.Range("D3").Value = '3.369.89'
$ws.Range("E3").Value = '  +4.62%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws "D5" '191.96'
$ws.Range("E5").Value = '  +5.21%  '
Set-TextValue $ws "D6" '593.60'
$ws.Range("E6").Value = '  +2.82%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws "D7" '0.610'
$ws.Range("E7").Value = '  +1.52%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws "D8" '1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +3.37%  '
Set-TextValue $ws "D10" '6.76'
$ws.Range("E10").Value = '  +3.64%  '
Set-TextValue $ws "D11" '0.422'
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("D12").Value = '3.958.77'
$ws.Range("E12").Value = '  +4.78%  '
$ws.Range("E13").Value = '  +1.24%  '
Set-TextValue $ws "D14" '28.70'
$ws.Range("E14").Value = '  +3.73%  '
$ws.Range("D15").Value = '69.557.21'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = '3.354.31'
$ws.Range("E17").Value = '  +4.08%  '
Set-TextValue $ws "D18" '449.03'
$ws.Range("E18").Value = '  +14.04%  '
Set-TextValue $ws "D19" '5.85'
$ws.Range("E19").Value = '  +1.92%  '
Set-TextValue $ws "D20" '13.83'
$ws.Range("E20").Value = '  +3.30%  '
Set-TextValue $ws "D21" '7.83'
$ws.Range("E21").Value = '  +3.85%  '
Set-TextValue $ws "D22" '73.58'
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '3.520.24'
$ws.Range("E24").Value = '  +4.69%  '
Set-TextValue $ws "D25" '0.520'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  +4.11%  '
$ws.Range("E27").Value = '  +4.68%  '
Set-TextValue $ws "D28" '9.61'
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("E29").Value = '  +0.03%  '
Set-TextValue $ws "D31" '23.26'
$ws.Range("E31").Value = '  +3.00%  '
Set-TextValue $ws "D32" '5.63'
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("E33").Value = '  +4.78%  '
Set-TextValue $ws "D34" '7.05'
$ws.Range("E34").Value = '  +1.72%  '
$ws.Range("E35").Value = '  +0.00%  '
Set-TextValue $ws "D36" '1.53'
$ws.Range("E36").Value = '  +4.21%  '
Set-TextValue $ws "D37" '165.05'
$ws.Range("E37").Value = '  +2.76%  '
$ws.Range("E38").Value = '  +3.63%  '
Set-TextValue $ws "D39" '27.33'
$ws.Range("E39").Value = '  +4.54%  '
Set-TextValue $ws "D40" '0.823'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +1.46%  '
Set-TextValue $ws "D42" '6.52'
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").Value = '2.746.39'
$ws.Range("E43").Value = '  +6.16%  '
Set-TextValue $ws "D44" '2.54'
$ws.Range("E44").Value = '  +3.36%  '
Set-TextValue $ws "D45" '25.67'
$ws.Range("E45").Value = '  +5.03%  '
Set-TextValue $ws "D46" '0.0690'
$ws.Range("E46").Value = '  +1.20%  '
Set-TextValue $ws "D47" '345.22'
$ws.Range("E47").Value = '  +3.82%  '
Set-TextValue $ws "D48" '40.81'
Set-TextValue $ws "D49" '0.0287'
$ws.Range("E49").Value = '  +3.82%  '
Set-TextValue $ws "D50" '1.03'
$ws.Range("E50").Value = '  +7.65%  '
Set-TextValue $ws "D51" '32.89'
$ws.Range("E51").Value = '  +7.55%  '
